$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 426-430: only the "Nb nouveaux cas positifs" input (column C) changes.
# Column B ("Cumul cas positifs") is a carry-forward formula (IF(TODAY()>...))
# and recalculates automatically once C is updated.
$ws.Range("C426").Value = 102
$ws.Range("C427").Value = 100
$ws.Range("C429").Value = 87
$ws.Range("C430").Value = 86

# Rows 431-433: these rows previously had no recorded data (inputs blank,
# carry-forward formulas resolving to ""). Fill in the newly published daily
# figures. Columns: C = new cases, E = SI patients, F = intubated,
# G = hospitalised outside SI, L = new hospital deaths,
# M = new extra-hospital deaths. H/J/K are formulas and recompute on their own.
#
# L/M are formatted as Text (numFmtId 49 / "@") in this sheet, and a plain
# .Value assignment on a Text-formatted cell is stored verbatim as text (same
# as typing into the cell). Flip the format to (lowercase) "general" - which
# resolves to the workbook's existing General+border style instead of minting
# a new one - write the numeric value, then flip back to "@" so the cell ends
# up numeric while keeping its original style index.
function Set-NumericOnTextCell($cell, $val) {
    $cell.NumberFormat = "general"
    $cell.Value = $val
    $cell.NumberFormat = "@"
}

$ws.Range("C431").Value = 42
$ws.Range("E431").Value = 7
$ws.Range("F431").Value = 6
$ws.Range("G431").Value = 27
Set-NumericOnTextCell $ws.Range("L431") 0
Set-NumericOnTextCell $ws.Range("M431") 0

$ws.Range("C432").Value = 12
$ws.Range("E432").Value = 8
$ws.Range("F432").Value = 7
$ws.Range("G432").Value = 28
Set-NumericOnTextCell $ws.Range("L432") 0
Set-NumericOnTextCell $ws.Range("M432") 0

$ws.Range("C433").Value = 7
$ws.Range("E433").Value = 7
$ws.Range("F433").Value = 6
$ws.Range("G433").Value = 33
Set-NumericOnTextCell $ws.Range("L433") 0
Set-NumericOnTextCell $ws.Range("M433") 0
